# Apply the daily cryptos-list price/volume refresh (GitHub Actions data pull).
# Row 35/36 additionally swap identity (ARBITRUM <-> MXToken reordered) per the feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.215.15"
$ws.Range("E2").Value = "  -0.23%  "

# Row 3 (Ethereum)
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.660.83"
$ws.Range("E3").Value = "  -0.14%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.46%  "

# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.43"
$ws.Range("E5").Value = "  +0.32%  "

# Row 6 (XRP)
$ws.Range("E6").Value = "  -0.79%  "

# Row 7 (USDC)
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("E7").Value = "  -0.43%  "

# Row 8 (Cardano)
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2637"
$ws.Range("E8").Value = "  +0.07%  "

# Row 9 (Dogecoin)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06317"
$ws.Range("E9").Value = "  -0.66%  "

# Row 10 (Solana)
$ws.Range("E10").Value = "  +0.60%  "

# Row 11 (TRON)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07804"
$ws.Range("E11").Value = "  -0.64%  "

# Row 12 (Polkadot)
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.500"
$ws.Range("E12").Value = "  -1.33%  "

# Row 13 (WrappedEther)
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.662.89"
$ws.Range("E13").Value = "  -0.07%  "

# Row 14 (WrappedliquidstakedEther2.0)
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.888.82"
$ws.Range("E14").Value = "  -0.14%  "

# Row 15 (Polygon)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5553"
$ws.Range("E15").Value = "  +0.57%  "

# Row 16 (ShibaInu)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8033"
$ws.Range("E16").Value = "  -1.69%  "

# Row 17 (Litecoin)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.27"
$ws.Range("E17").Value = "  -0.54%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.218.88"
$ws.Range("E18").Value = "  -0.31%  "

# Row 19 (Dai)
$ws.Range("E19").Value = "  -0.41%  "

# Row 20 (Uniswap)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.645"
$ws.Range("E20").Value = "  -0.53%  "

# Row 21 (BitcoinCash)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "196.80"
$ws.Range("E21").Value = "  +2.10%  "

# Row 22 (Avalanche)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.14"
$ws.Range("E22").Value = "  -0.68%  "

# Row 23 (Chainlink)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.979"
$ws.Range("E23").Value = "  -0.82%  "

# Row 24 (BinanceUSD)
$ws.Range("E24").Value = "  -0.50%  "

# Row 25 (Monero)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.75"
$ws.Range("E25").Value = "  +1.01%  "

# Row 26 (Stellar)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1206"
$ws.Range("E26").Value = "  -1.59%  "

# Row 27 (Cosmos)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.154"
$ws.Range("E27").Value = "  -0.54%  "

# Row 28 (EthereumClassic)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.06"
$ws.Range("E28").Value = "  +0.00%  "

# Row 29 (Toncoin)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.503"
$ws.Range("E29").Value = "  +1.67%  "

# Row 30 (Hedera)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05779"
$ws.Range("E30").Value = "  -2.34%  "

# Row 31 (PancakeSwap)
$ws.Range("E31").Value = "  -0.42%  "

# Row 32 (InternetComputer(DFINITY))
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.492"
$ws.Range("E32").Value = "  -2.64%  "

# Row 33 (Filecoin)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.353"
$ws.Range("E33").Value = "  +2.31%  "

# Row 34 (LidoDAOToken)
$ws.Range("E34").Value = "  -1.78%  "

# Row 35 (MXToken)
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9555"
$ws.Range("E35").Value = "  -0.38%  "

# Row 36 (ARBITRUM)
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.810"
$ws.Range("E36").Value = "  -0.62%  "

# Row 37 (HuobiToken)
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.421"
$ws.Range("E37").Value = "  -0.20%  "

# Row 38 (ImmutableX)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5737"
$ws.Range("E38").Value = "  -0.84%  "

# Row 39 (VeChain)
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01595"
$ws.Range("E39").Value = "  -0.38%  "

# Row 40 (FraxShare)
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.965"
$ws.Range("E40").Value = "  +1.90%  "

# Row 41 (Maker)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.063.09"
$ws.Range("E41").Value = "  +1.70%  "

# Row 42 (TrustWalletToken)
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8512"
$ws.Range("E42").Value = "  -1.57%  "

# Row 43 (PaxDollar)
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.006"
$ws.Range("E43").Value = "  -0.39%  "

# Row 44 (Quant)
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.29"
$ws.Range("E44").Value = "  -0.65%  "

# Row 45 (RocketPoolETH)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.799.72"
$ws.Range("E45").Value = "  -0.25%  "

# Row 46 (Aave)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.29"
$ws.Range("E46").Value = "  +1.70%  "

# Row 47 (Frax)
$ws.Range("E47").Value = "  -0.14%  "

# Row 48 (Mantle)
$ws.Range("E48").Value = "  +0.57%  "

# Row 49 (EnergySwap)
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.036"
$ws.Range("E49").Value = "  +0.16%  "

# Row 50 (Cronos)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05206"
$ws.Range("E50").Value = "  +0.80%  "

# Row 51 (BabyDogeCoin)
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₈100"
$ws.Range("E51").Value = "  -5.96%  "
